$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Version bump
$meta.Range("B3").Value = "6.0.0"

# Date bump
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a value
$meta.Range("B9").Value = "Alvearie Team"

# "Contact" / "No display for ContactDetail" row becomes "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Drop the now-redundant second "Contact" row entirely (rows below shift up)
$meta.Rows.Item(11).Delete()

# --- Elements sheet ---
$elements = $wb.Worksheets.Item("Elements")

# Root Extension row: Short / Definition now describe this specific extension
$elements.Range("K2").Value = "Episode Procedure Type Code"
$elements.Range("L2").Value = "Primary procedure type AHRQ code for the episode of care"
